# Edit: iTalent work item 79 - English presentation
# 1) Slide 2 ("I-Talent Platform" bullet list): capitalize / re-run the
#    bullet lines "web application", "kickstarter", "pitch idea with media",
#    "like/subscribe" and "status updates" without touching the <a:br/>
#    soft line breaks already in the title text box.
# 2) Slide 3: move/resize the "Content Placeholder 2" shape and rename the
#    last "Documentation" bullet to "Quality Management".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 - split each "- <phrase>" run into a literal "- " run plus a
# capitalized run (and, for the "pitch idea with media" line, a third run),
# using precise character offsets so the existing <a:br/> breaks and run
# formatting (rPr) are left completely untouched.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange

$tr2.Characters(20, 2).Text  = "- "
$tr2.Characters(22, 15).Text = "Web Application"

$tr2.Characters(38, 2).Text  = "- "
$tr2.Characters(40, 11).Text = "Kickstarter"

$tr2.Characters(52, 2).Text  = "- "
$tr2.Characters(54, 6).Text  = "Pitch "
$tr2.Characters(60, 15).Text = "idea with media"

$tr2.Characters(76, 2).Text  = "- "
$tr2.Characters(78, 14).Text = "Like/Subscribe"

$tr2.Characters(93, 2).Text  = "- "
$tr2.Characters(95, 14).Text = "Status Updates"

# ---------------------------------------------------------------------
# Slide 3 - reposition/resize the content placeholder and rewrite the
# final "Documentation" paragraph to "Quality Management".
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$contentShape = $s3.Shapes.Item(2)

$contentShape.Left   = 206.9121322631836
$contentShape.Top    = 139.59180450439456
$contentShape.Width  = 317.51557922363287
$contentShape.Height = 356.3750457763672

$tr3 = $contentShape.TextFrame.TextRange
$tr3.Paragraphs(7).Text = "Quality Management"
